$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H107").Value = 167633
$ws.Range("I107").Value = 250949.75
$ws.Range("J107").Value = 999.5
$ws.Range("K107").Value = 250949.75
$ws.Range("L107").Value = 999.5
$ws.Range("M107").Value = -249029.75
$ws.Range("N107").Value = -4839.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 7350.3477
$ws.Range("J112").Value = 8322.9
$ws.Range("L112").Value = 24968.7
$ws.Range("N112").Value = -27184.7

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 4654299
$ws.Range("I116").Value = 5408470.5
$ws.Range("J116").Value = 3574.5
$ws.Range("K116").Value = 5408470.5
$ws.Range("L116").Value = 3574.5
$ws.Range("M116").Value = -5405028.5
$ws.Range("N116").Value = -10458.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H129").Value = 1133.8667
$ws.Range("J129").Value = 1144.7441
$ws.Range("L129").Value = 3434.2323
$ws.Range("N129").Value = -13434.2323

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 1859.9636
$ws.Range("I132").Value = 1743.48
$ws.Range("J132").Value = 3024.8
$ws.Range("K132").Value = 5230.440000000001
$ws.Range("L132").Value = 9074.400000000001
$ws.Range("M132").Value = -2700.440000000001
$ws.Range("N132").Value = -14134.4

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H134").Value = 117006.31
$ws.Range("J134").Value = 117006.31
$ws.Range("L134").Value = 117006.31
$ws.Range("N134").Value = -127146.31

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 1897.0294
$ws.Range("I137").Value = 1464.1786
$ws.Range("K137").Value = 4392.5358
$ws.Range("M137").Value = -1842.5358

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 8342131.5
$ws.Range("I138").Value = 15389244
$ws.Range("J138").Value = 13725.909
$ws.Range("K138").Value = 46167732
$ws.Range("L138").Value = 41177.727
$ws.Range("M138").Value = -46162592
$ws.Range("N138").Value = -51457.727

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H139").Value = 72219.164
$ws.Range("J139").Value = 72219.164
$ws.Range("L139").Value = 72219.164
$ws.Range("N139").Value = -82499.164

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 51528.832
$ws.Range("I32").Value = 41071.23
$ws.Range("K32").Value = 41071.23
$ws.Range("M32").Value = -40784.23

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 72698.5
$ws.Range("I102").Value = 1136.5555
$ws.Range("K102").Value = 1136.5555
$ws.Range("M102").Value = 485.4445000000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H134").Value = 55555
$ws.Range("J134").Value = 55555
$ws.Range("L134").Value = 55555
$ws.Range("N134").Value = -65695

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H135").Value = 57000
$ws.Range("J135").Value = 57000
$ws.Range("L135").Value = 57000
$ws.Range("N135").Value = -67140

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H132").Value = 79561.89
$ws.Range("J132").Value = 79561.89
$ws.Range("L132").Value = 79561.89
$ws.Range("N132").Value = -89681.89

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5366.4287
$ws.Range("I31").Value = 4613.6
$ws.Range("K31").Value = 4613.6
$ws.Range("M31").Value = -4318.6

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 5366.4287
$ws.Range("I34").Value = 4613.6
$ws.Range("K34").Value = 4613.6
$ws.Range("M34").Value = -4411.6

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H104").Value = 0
$ws.Range("J104").Value = 0
$ws.Range("L104").Value = 0
$ws.Range("N104").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 5563.273
$ws.Range("I105").Value = 6493.647
$ws.Range("K105").Value = 6493.647
$ws.Range("M105").Value = -4746.647

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 1803.973
$ws.Range("I132").Value = 1621.3715
$ws.Range("J132").Value = 4999.5
$ws.Range("K132").Value = 4864.1145
$ws.Range("L132").Value = 14998.5
$ws.Range("M132").Value = -2334.1145
$ws.Range("N132").Value = -20058.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 1268.8837
$ws.Range("I134").Value = 1208.7
$ws.Range("J134").Value = 2071.3333
$ws.Range("K134").Value = 3626.1
$ws.Range("L134").Value = 6213.999899999999
$ws.Range("M134").Value = -1091.1
$ws.Range("N134").Value = -11283.9999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H135").Value = 139779.8
$ws.Range("J135").Value = 139779.8
$ws.Range("L135").Value = 139779.8
$ws.Range("N135").Value = -149919.8

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H138").Value = 82846
$ws.Range("J138").Value = 82846
$ws.Range("L138").Value = 82846
$ws.Range("N138").Value = -93126

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H140").Value = 65333.332
$ws.Range("J140").Value = 78000
$ws.Range("L140").Value = 78000
$ws.Range("N140").Value = -88360

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H75").Value = 4911.3335
$ws.Range("J75").Value = 9750
$ws.Range("L75").Value = 29250
$ws.Range("N75").Value = -31246

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H78").Value = 4911.3335
$ws.Range("J78").Value = 9750
$ws.Range("L78").Value = 87750
$ws.Range("N78").Value = -97734

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 11653.5
$ws.Range("I113").Value = 908
$ws.Range("J113").Value = 18817.166
$ws.Range("K113").Value = 2724
$ws.Range("L113").Value = 56451.49800000001
$ws.Range("M113").Value = -554
$ws.Range("N113").Value = -60791.49800000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H130").Value = 2600
$ws.Range("I130").Value = 1000
$ws.Range("K130").Value = 3000
$ws.Range("M130").Value = 2020

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 33338546
$ws.Range("J131").Value = 37039068
$ws.Range("L131").Value = 111117204
$ws.Range("N131").Value = -111127284

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H123").Value = 14531.25
$ws.Range("J123").Value = 14531.25
$ws.Range("L123").Value = 14531.25
$ws.Range("N123").Value = -19431.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H133").Value = 54033.332
$ws.Range("J133").Value = 54033.332
$ws.Range("L133").Value = 54033.332
$ws.Range("N133").Value = -64153.332

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4536.6665
$ws.Range("I7").Value = 4495.4546
$ws.Range("K7").Value = 4495.4546
$ws.Range("M7").Value = -4383.4546

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H24").Value = 13498
$ws.Range("J24").Value = 13498
$ws.Range("L24").Value = 13498
$ws.Range("N24").Value = -14184

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 4536.6665
$ws.Range("I126").Value = 4495.4546
$ws.Range("K126").Value = 13486.3638
$ws.Range("M126").Value = -11016.3638

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H133").Value = 70995.375
$ws.Range("J133").Value = 70995.375
$ws.Range("L133").Value = 70995.375
$ws.Range("N133").Value = -76055.375

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H134").Value = 39880
$ws.Range("J134").Value = 39880
$ws.Range("L134").Value = 39880
$ws.Range("N134").Value = -50020

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H140").Value = 84750
$ws.Range("J140").Value = 84750
$ws.Range("L140").Value = 84750
$ws.Range("N140").Value = -95110

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1387.2449
$ws.Range("I132").Value = 1259.5555
$ws.Range("J132").Value = 1740.8462
$ws.Range("K132").Value = 3778.6665
$ws.Range("L132").Value = 5222.5386
$ws.Range("M132").Value = -1248.6665
$ws.Range("N132").Value = -10282.5386

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H139").Value = 62835.883
$ws.Range("J139").Value = 62835.883
$ws.Range("L139").Value = 62835.883
$ws.Range("N139").Value = -73115.883

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H141").Value = 69440
$ws.Range("J141").Value = 69440
$ws.Range("L141").Value = 69440
$ws.Range("N141").Value = -79800
